$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 66687024
$ws.Range("I62").Value = 125003450
$ws.Range("J62").Value = 39685.57
$ws.Range("K62").Value = 125003450
$ws.Range("L62").Value = 39685.57
$ws.Range("M62").Value = -125002826
$ws.Range("N62").Value = -40933.57

$ws.Range("H65").Value = 66687024
$ws.Range("I65").Value = 125003450
$ws.Range("J65").Value = 39685.57
$ws.Range("K65").Value = 625017250
$ws.Range("L65").Value = 198427.85
$ws.Range("M65").Value = -625014130
$ws.Range("N65").Value = -204667.85

$ws.Range("H76").Value = 19808.75
$ws.Range("I76").Value = 19808.75
$ws.Range("K76").Value = 19808.75
$ws.Range("M76").Value = -19493.75

$ws.Range("H79").Value = 19808.75
$ws.Range("I79").Value = 19808.75
$ws.Range("K79").Value = 19808.75
$ws.Range("M79").Value = -18716.75

$ws.Range("H92").Value = 5350
$ws.Range("I92").Value = 699.5
$ws.Range("J92").Value = 10000.5
$ws.Range("K92").Value = 699.5
$ws.Range("L92").Value = 10000.5
$ws.Range("M92").Value = 548.5
$ws.Range("N92").Value = -12496.5

$ws.Range("H99").Value = 495.5
$ws.Range("I99").Value = 495.5
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1486.5
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 11.5
$ws.Range("N99").ClearContents()

$ws.Range("H111").Value = 15632114
$ws.Range("I111").Value = 17864278
$ws.Range("K111").Value = 53592834
$ws.Range("M111").Value = -53589767

$ws.Range("H116").Value = 14711554
$ws.Range("I116").Value = 62502660
$ws.Range("J116").Value = 6599
$ws.Range("K116").Value = 62502660
$ws.Range("L116").Value = 6599
$ws.Range("M116").Value = -62499218
$ws.Range("N116").Value = -13483

$ws.Range("H125").Value = 31251674
$ws.Range("J125").Value = 3996
$ws.Range("L125").Value = 35964
$ws.Range("N125").Value = -40884

$ws.Range("H138").Value = 5857.51
$ws.Range("J138").Value = 7334.7837
$ws.Range("L138").Value = 22004.3511
$ws.Range("N138").Value = -32284.3511

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 1848
$ws.Range("I16").Value = 395
$ws.Range("J16").Value = 2816.6667
$ws.Range("K16").Value = 395
$ws.Range("L16").Value = 2816.6667
$ws.Range("M16").Value = -108
$ws.Range("N16").Value = -3390.6667

$ws.Range("H61").Value = 9249.200000000001
$ws.Range("I61").Value = 2454.3333
$ws.Range("J61").Value = 14808.637
$ws.Range("K61").Value = 2454.3333
$ws.Range("L61").Value = 14808.637
$ws.Range("M61").Value = -2242.3333
$ws.Range("N61").Value = -15232.637

$ws.Range("H110").Value = 23810448
$ws.Range("I110").Value = 972.7
$ws.Range("J110").Value = 83334136
$ws.Range("K110").Value = 972.7
$ws.Range("L110").Value = 83334136
$ws.Range("M110").Value = 1072.3
$ws.Range("N110").Value = -83338226

$ws.Range("H136").Value = 9249.200000000001
$ws.Range("I136").Value = 2454.3333
$ws.Range("J136").Value = 14808.637
$ws.Range("K136").Value = 7362.999899999999
$ws.Range("L136").Value = 44425.911
$ws.Range("M136").Value = -4812.999899999999
$ws.Range("N136").Value = -49525.911

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 59376
$ws.Range("J55").Value = 59376
$ws.Range("L55").Value = 59376
$ws.Range("N55").Value = -59922

$ws.Range("H94").Value = 1110.4
$ws.Range("I94").Value = 700.5714
$ws.Range("J94").Value = 2066.6667
$ws.Range("K94").Value = 700.5714
$ws.Range("L94").Value = 2066.6667
$ws.Range("M94").Value = -249.5714
$ws.Range("N94").Value = -2968.6667

$ws.Range("H99").Value = 6994378
$ws.Range("I99").Value = 858
$ws.Range("K99").Value = 858
$ws.Range("M99").Value = 640

$ws.Range("H134").Value = 5042.5093
$ws.Range("I134").Value = 2416.457
$ws.Range("K134").Value = 7249.370999999999
$ws.Range("M134").Value = -4714.370999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 54188
$ws.Range("J18").Value = 54188
$ws.Range("L18").Value = 54188
$ws.Range("N18").Value = -54648

$ws.Range("H31").Value = 6614.262
$ws.Range("I31").Value = 2811.0278
$ws.Range("J31").Value = 12090.92
$ws.Range("K31").Value = 2811.0278
$ws.Range("L31").Value = 12090.92
$ws.Range("M31").Value = -2516.0278
$ws.Range("N31").Value = -12680.92

$ws.Range("H34").Value = 6614.262
$ws.Range("I34").Value = 2811.0278
$ws.Range("J34").Value = 12090.92
$ws.Range("K34").Value = 2811.0278
$ws.Range("L34").Value = 12090.92
$ws.Range("M34").Value = -2609.0278
$ws.Range("N34").Value = -12494.92

$ws.Range("H59").Value = 97496
$ws.Range("J59").Value = 97496
$ws.Range("L59").Value = 97496
$ws.Range("N59").Value = -99786

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 5900
$ws.Range("I62").Value = 5000
$ws.Range("K62").Value = 15000
$ws.Range("M62").Value = -14314

$ws.Range("H65").Value = 5900
$ws.Range("I65").Value = 5000
$ws.Range("K65").Value = 45000
$ws.Range("M65").Value = -41568

$ws.Range("H87").Value = 71438700
$ws.Range("J87").Value = 12000
$ws.Range("L87").Value = 36000
$ws.Range("N87").Value = -38496

$ws.Range("H90").Value = 71438700
$ws.Range("J90").Value = 12000
$ws.Range("L90").Value = 108000
$ws.Range("N90").Value = -120480

$ws.Range("H121").Value = 25000646
$ws.Range("I121").Value = 50000170
$ws.Range("J121").Value = 16667471
$ws.Range("K121").Value = 150000510
$ws.Range("L121").Value = 50002413
$ws.Range("M121").Value = -149999200
$ws.Range("N121").Value = -50005033

$ws.Range("H132").Value = 10660.059
$ws.Range("J132").Value = 13150.8
$ws.Range("L132").Value = 118357.2
$ws.Range("N132").Value = -123417.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1538589.2
$ws.Range("I2").Value = 74
$ws.Range("J2").Value = 2857316.8
$ws.Range("K2").Value = 74
$ws.Range("L2").Value = 2857316.8
$ws.Range("M2").Value = 39
$ws.Range("N2").Value = -2857542.8

$ws.Range("H97").Value = 1596.1052
$ws.Range("I97").Value = 1397.5
$ws.Range("K97").Value = 1397.5
$ws.Range("M97").Value = -901.5

$ws.Range("H98").Value = 65068.4
$ws.Range("J98").Value = 65068.4
$ws.Range("L98").Value = 65068.4
$ws.Range("N98").Value = -71058.39999999999

$ws.Range("H113").Value = 8478.130999999999
$ws.Range("I113").Value = 6250
$ws.Range("J113").Value = 8947.210999999999
$ws.Range("K113").Value = 6250
$ws.Range("L113").Value = 8947.210999999999
$ws.Range("M113").Value = -4080
$ws.Range("N113").Value = -13287.211

$ws.Range("H132").Value = 7716.727
$ws.Range("I132").Value = 3110.5
$ws.Range("K132").Value = 9331.5
$ws.Range("M132").Value = -6801.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()

$ws.Range("H22").Value = 12599.8
$ws.Range("I22").Value = 4226.909
$ws.Range("J22").Value = 22833.334
$ws.Range("K22").Value = 4226.909
$ws.Range("L22").Value = 22833.334
$ws.Range("M22").Value = -3931.909
$ws.Range("N22").Value = -23423.334

$ws.Range("H27").Value = 12599.8
$ws.Range("I27").Value = 4226.909
$ws.Range("J27").Value = 22833.334
$ws.Range("K27").Value = 4226.909
$ws.Range("L27").Value = 22833.334
$ws.Range("M27").Value = -4119.909
$ws.Range("N27").Value = -23047.334

$ws.Range("H40").Value = 4641.353
$ws.Range("I40").Value = 3709.8
$ws.Range("J40").Value = 5972.143
$ws.Range("K40").Value = 3709.8
$ws.Range("L40").Value = 5972.143
$ws.Range("M40").Value = -3573.8
$ws.Range("N40").Value = -6244.143

$ws.Range("H82").Value = 2645.8333
$ws.Range("I82").Value = 2488.4546
$ws.Range("J82").Value = 2893.1428
$ws.Range("K82").Value = 2488.4546
$ws.Range("L82").Value = 2893.1428
$ws.Range("M82").Value = -2127.4546
$ws.Range("N82").Value = -3615.1428

$ws.Range("H85").Value = 2645.8333
$ws.Range("I85").Value = 2488.4546
$ws.Range("J85").Value = 2893.1428
$ws.Range("K85").Value = 2488.4546
$ws.Range("L85").Value = 2893.1428
$ws.Range("M85").Value = -1240.4546
$ws.Range("N85").Value = -5389.1428

$ws.Range("H93").Value = 4047.4707
$ws.Range("I93").Value = 4680.7
$ws.Range("J93").Value = 3142.8572
$ws.Range("K93").Value = 4680.7
$ws.Range("L93").Value = 3142.8572
$ws.Range("M93").Value = -3432.7
$ws.Range("N93").Value = -5638.8572

$ws.Range("H100").Value = 3021.1
$ws.Range("J100").Value = 4002.6667
$ws.Range("L100").Value = 4002.6667
$ws.Range("N100").Value = -5084.6667

$ws.Range("H122").Value = 3554.7046
$ws.Range("I122").Value = 2803.2122
$ws.Range("J122").Value = 5809.1816
$ws.Range("K122").Value = 8409.6366
$ws.Range("L122").Value = 17427.5448
$ws.Range("M122").Value = -5959.6366
$ws.Range("N122").Value = -22327.5448

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 20842662
$ws.Range("I132").Value = 45459908
$ws.Range("K132").Value = 136379724
$ws.Range("M132").Value = -136377194

$ws.Range("H136").Value = 9737.861999999999
$ws.Range("I136").Value = 4866.0835
$ws.Range("K136").Value = 14598.2505
$ws.Range("M136").Value = -12048.2505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1710.7142
$ws.Range("I96").Value = 1579.3334
$ws.Range("K96").Value = 1579.3334
$ws.Range("M96").Value = -206.3334
